$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values
$ws.Range('D2').Value = '35.306.60'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '1.898.38'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.692'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.10%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.35'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.349'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.24%  '
$ws.Range('E11').Value = '  +2.30%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = '2.172.60'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.708'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '1.896.61'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '35.251.81'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('D20').Value = '0.0₃0820'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  +1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.29%  '
$ws.Range('E30').Value = '  +3.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0569'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.73%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.911'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.10%  '
$ws.Range('E38').Value = '  +7.37%  '
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '95.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0656'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0207'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').Value = '1.356.00'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.49%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.06%  '
$ws.Range('E51').Value = '  -2.80%  '
